$wb = $excel.ActiveWorkbook

# --- Section_A timetable updates ---
$ws1 = $wb.Worksheets.Item("Section_A")
$ws1.Range("B2").Value = "ELECTIVE_B1 [C305]"
$ws1.Range("C2").Value = "DS161 [C104]"
$ws1.Range("D2").Value = "ELECTIVE_B1 [C305]"
$ws1.Range("E2").Value = "Free"
$ws1.Range("F2").Value = "Free"
$ws1.Range("B3").Value = "EC161 [C002]"
$ws1.Range("C3").Value = "HS161 [C204]"
$ws1.Range("D3").Value = "MA161 [C204]"
$ws1.Range("E3").Value = "DS161 [C104]"
$ws1.Range("F3").Value = "MA161 [C204]"
$ws1.Range("B5").Value = "HS161 [C204]"
$ws1.Range("C5").Value = "EC161 [C002]"
$ws1.Range("D5").Value = "Free"
$ws1.Range("E5").Value = "MA162 [C205]"
$ws1.Range("F5").Value = "EC161 (Lab) [L407]"
$ws1.Range("B6").Value = "Free"
$ws1.Range("C6").Value = "Free"
$ws1.Range("D6").Value = "Free"
$ws1.Range("E6").Value = "Free"
$ws1.Range("F6").Value = "EC161 (Lab) [L407]"
$ws1.Range("B7").Value = "Free"
$ws1.Range("C7").Value = "Free"
$ws1.Range("D7").Value = "Free"
$ws1.Range("E7").Value = "Free"
$ws1.Range("F7").Value = "MA162 [C205]"

# --- Section_B timetable updates ---
$ws2 = $wb.Worksheets.Item("Section_B")
$ws2.Range("B2").Value = "ELECTIVE_B1 [C402]"
$ws2.Range("C2").Value = "HS161 [C305]"
$ws2.Range("D2").Value = "ELECTIVE_B1 [C402]"
$ws2.Range("E2").Value = "Free"
$ws2.Range("F2").Value = "MA162 [C402]"
$ws2.Range("B3").Value = "MA161 [C302]"
$ws2.Range("C3").Value = "Free"
$ws2.Range("D3").Value = "MA161 [C302]"
$ws2.Range("E3").Value = "DS161 [C202]"
$ws2.Range("F3").Value = "Free"
$ws2.Range("B5").Value = "Free"
$ws2.Range("C5").Value = "Free"
$ws2.Range("D5").Value = "HS161 [C305]"
$ws2.Range("E5").Value = "EC161 [C203]"
$ws2.Range("F5").Value = "DS161 [C202]"
$ws2.Range("B6").Value = "Free"
$ws2.Range("C6").Value = "Free"
$ws2.Range("D6").Value = "Free"
$ws2.Range("E6").Value = "Free"
$ws2.Range("F6").Value = "Free"
$ws2.Range("B7").Value = "Free"
$ws2.Range("C7").Value = "Free"
$ws2.Range("D7").Value = "EC161 [C203]"
$ws2.Range("E7").Value = "MA162 [C402]"
$ws2.Range("F7").Value = "EC161 (Lab) [L306]"
$ws2.Range("B8").Value = "Free"
$ws2.Range("C8").Value = "Free"
$ws2.Range("D8").Value = "Free"
$ws2.Range("E8").Value = "Free"
$ws2.Range("F8").Value = "EC161 (Lab) [L306]"

# --- Classroom_Utilization updates (numeric Weekly/Daily/Utilization stats) ---
$ws3 = $wb.Worksheets.Item("Classroom_Utilization")
$ws3.Range("D3").Value = 3
$ws3.Range("E3").Value = 0.6
$ws3.Range("G3").Value = 7.5
$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 0
$ws3.Range("G4").Value = 0
$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 0
$ws3.Range("G5").Value = 0
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 0
$ws3.Range("G6").Value = 0
$ws3.Range("D13").Value = 0
$ws3.Range("E13").Value = 0
$ws3.Range("G13").Value = 0
$ws3.Range("D14").Value = 3
$ws3.Range("E14").Value = 0.6
$ws3.Range("G14").Value = 7.5
$ws3.Range("D15").Value = 3
$ws3.Range("E15").Value = 0.6
$ws3.Range("G15").Value = 7.5
$ws3.Range("D16").Value = 6
$ws3.Range("E16").Value = 1.2
$ws3.Range("G16").Value = 15
$ws3.Range("D22").Value = 3
$ws3.Range("E22").Value = 0.6
$ws3.Range("G22").Value = 7.5
$ws3.Range("D25").Value = 6
$ws3.Range("E25").Value = 1.2
$ws3.Range("G25").Value = 15
$ws3.Range("D29").Value = 0
$ws3.Range("E29").Value = 0
$ws3.Range("G29").Value = 0
$ws3.Range("D30").Value = 6
$ws3.Range("E30").Value = 1.2
$ws3.Range("G30").Value = 15
$ws3.Range("D33").Value = 0
$ws3.Range("E33").Value = 0
$ws3.Range("G33").Value = 0
$ws3.Range("D35").Value = 2.5
$ws3.Range("E35").Value = 0.5
$ws3.Range("G35").Value = 6.25
$ws3.Range("D36").Value = 0
$ws3.Range("E36").Value = 0
$ws3.Range("G36").Value = 0

# --- Classroom_Allocation updates ---
$ws4 = $wb.Worksheets.Item("Classroom_Allocation")
$ws4.Range("G2").Value = "C305"
$ws4.Range("H2").Value = "classroom"
$ws4.Range("I2").Value = "'96"
$ws4.Range("G3").Value = "C002"
$ws4.Range("H3").Value = "large classroom"
$ws4.Range("I3").Value = "'116"
$ws4.Range("F4").Value = "HS161"
$ws4.Range("G4").Value = "C204"
$ws4.Range("I4").Value = "'96"
$ws4.Range("E5").Value = "09:00-10:30"
$ws4.Range("F5").Value = "DS161"
$ws4.Range("E6").Value = "10:30-12:00"
$ws4.Range("F6").Value = "HS161"
$ws4.Range("G6").Value = "C204"
$ws4.Range("H6").Value = "classroom"
$ws4.Range("I6").Value = "'96"
$ws4.Range("E7").Value = "13:00-14:30"
$ws4.Range("F7").Value = "EC161"
$ws4.Range("G7").Value = "C002"
$ws4.Range("H7").Value = "large classroom"
$ws4.Range("I7").Value = "'116"
$ws4.Range("G8").Value = "C305"
$ws4.Range("H8").Value = "classroom"
$ws4.Range("I8").Value = "'96"
$ws4.Range("F9").Value = "MA161"
$ws4.Range("G9").Value = "C204"
$ws4.Range("H9").Value = "classroom"
$ws4.Range("I9").Value = "'96"
$ws4.Range("D10").Value = "Thu"
$ws4.Range("E10").Value = "10:30-12:00"
$ws4.Range("F10").Value = "DS161"
$ws4.Range("E11").Value = "13:00-14:30"
$ws4.Range("F11").Value = "MA162"
$ws4.Range("G11").Value = "C205"
$ws4.Range("D12").Value = "Fri"
$ws4.Range("E12").Value = "10:30-12:00"
$ws4.Range("F12").Value = "MA161"
$ws4.Range("G12").Value = "C204"
$ws4.Range("H12").Value = "classroom"
$ws4.Range("I12").Value = "'96"
$ws4.Range("J12").Value = "Projector"
$ws4.Range("D13").Value = "Fri"
$ws4.Range("E13").Value = "13:00-14:30"
$ws4.Range("G13").Value = "L407"
$ws4.Range("H13").Value = "classroom"
$ws4.Range("D14").Value = "Fri"
$ws4.Range("E14").Value = "14:30-15:30"
$ws4.Range("F14").Value = "EC161 (Lab)"
$ws4.Range("G14").Value = "L407"
$ws4.Range("J14").Value = "Computers"
$ws4.Range("E15").Value = "15:30-17:00"
$ws4.Range("F15").Value = "MA162"
$ws4.Range("G15").Value = "C205"
$ws4.Range("I15").Value = "'96"
$ws4.Range("G16").Value = "C402"
$ws4.Range("E17").Value = "10:30-12:00"
$ws4.Range("F17").Value = "MA161"
$ws4.Range("G17").Value = "C302"
$ws4.Range("D18").Value = "Tue"
$ws4.Range("E18").Value = "09:00-10:30"
$ws4.Range("G18").Value = "C305"
$ws4.Range("G19").Value = "C402"
$ws4.Range("G20").Value = "C302"
$ws4.Range("F21").Value = "HS161"
$ws4.Range("G21").Value = "C305"
$ws4.Range("H21").Value = "classroom"
$ws4.Range("I21").Value = "'96"
$ws4.Range("F22").Value = "EC161"
$ws4.Range("G22").Value = "C203"
$ws4.Range("H22").Value = "classroom"
$ws4.Range("I22").Value = "'96"
$ws4.Range("E23").Value = "10:30-12:00"
$ws4.Range("F23").Value = "DS161"
$ws4.Range("G23").Value = "C202"
$ws4.Range("F24").Value = "EC161"
$ws4.Range("G24").Value = "C203"
$ws4.Range("F25").Value = "MA162"
$ws4.Range("G25").Value = "C402"
$ws4.Range("H25").Value = "classroom"
$ws4.Range("I25").Value = "'96"
$ws4.Range("E26").Value = "09:00-10:30"
$ws4.Range("F26").Value = "MA162"
$ws4.Range("G26").Value = "C402"
$ws4.Range("F27").Value = "DS161"
$ws4.Range("G27").Value = "C202"
$ws4.Range("J27").Value = "Projector"
$ws4.Range("E28").Value = "15:30-17:00"
$ws4.Range("E29").Value = "17:00-18:00"
$ws4.Range("F29").Value = "EC161 (Lab)"
$ws4.Range("G29").Value = "L306"
$ws4.Range("H29").Value = "classroom"
$ws4.Range("I29").Value = "'96"
$ws4.Range("J29").Value = "Computers"
